$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.211.70"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "1.809.21"
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.28"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -2.28%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3564"
$c.ClearFormats()
$ws.Range("E8").Value = "  -4.46%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07134"
$c.ClearFormats()
$ws.Range("E9").Value = "  -3.84%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8501"
$c.ClearFormats()
$ws.Range("E10").Value = "  -3.93%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.18"
$c.ClearFormats()
$ws.Range("D12").Value = "1.800.27"
$ws.Range("E12").Value = "  -10.35%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.318"
$c.ClearFormats()
$ws.Range("E13").Value = "  -3.28%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.366"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.08%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06849"
$c.ClearFormats()
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("E16").Value = "  -0.05%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "80.93"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.49%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008773"
$c.ClearFormats()
$ws.Range("E18").Value = "  -4.00%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").Value = "27.237.36"
$ws.Range("E21").Value = "  -3.77%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.106"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.40%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.ClearFormats()
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "2.047.58"
$ws.Range("E24").Value = "  -7.98%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.970"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.01%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "153.79"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.38%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.17"
$c.ClearFormats()
$ws.Range("E27").Value = "  -3.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.083"
$c.ClearFormats()
$ws.Range("E28").Value = "  -6.04%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "113.47"
$c.ClearFormats()
$ws.Range("E29").Value = "  -3.40%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.679"
$c.ClearFormats()
$ws.Range("E30").Value = "  -10.10%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08907"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.01%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7378"
$c.ClearFormats()
$ws.Range("E32").Value = "  -6.87%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.937"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.45%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.433"
$c.ClearFormats()
$ws.Range("E34").Value = "  -5.58%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.106"
$c.ClearFormats()
$ws.Range("E35").Value = "  -6.53%  "
$ws.Range("E36").Value = "  -0.04%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.075"
$c.ClearFormats()
$ws.Range("E37").Value = "  -4.86%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05183"
$c.ClearFormats()
$ws.Range("E38").Value = "  -5.10%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01905"
$c.ClearFormats()
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.721"
$c.ClearFormats()
$ws.Range("E40").Value = "  -5.90%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1633"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.34%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4971"
$c.ClearFormats()
$ws.Range("E42").Value = "  -3.93%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.270"
$c.ClearFormats()
$ws.Range("E43").Value = "  -9.11%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.178"
$c.ClearFormats()
$ws.Range("E44").Value = "  -5.79%  "
$ws.Range("E45").Value = "  -1.12%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.ClearFormats()
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  -0.09%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.06364"
$c.ClearFormats()
$ws.Range("E48").Value = "  -3.68%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.4565"
$c.ClearFormats()
$ws.Range("E49").Value = "  -4.20%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.596"
$c.ClearFormats()
$ws.Range("E50").Value = "  -3.71%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "62.56"
$c.ClearFormats()
$ws.Range("E51").Value = "  -5.01%  "
